# Core : add print qrcode
# Adds a new "sumber_dana" column (F) with value "Dana BOS" for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1: same bold header style as the rest of row 1 ---
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "sumber_dana"

# --- Data cells F2:F4: same body style as column A/B (vertical-centered) ---
$ws.Range("A2").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122)
$ws.Range("F2").Value = "Dana BOS"
$ws.Range("F3").Value = "Dana BOS"
$ws.Range("F4").Value = "Dana BOS"

$excel.CutCopyMode = 0

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 12.67

# --- Selection mirrors the edited column ---
$ws.Range("F2:F4").Select() | Out-Null
